$wb = $excel.ActiveWorkbook

# --- "Play Solitaire" sheet: fill in the remaining Main Success Scenario steps ---
$wsPlay = $wb.Worksheets.Item("Play Solitaire")
$wsPlay.Range("C12").Value = "Solitaire game"
$wsPlay.Range("B13").Value = "wins or loses solitaire game"
$wsPlay.Range("C13").Value = "displays return to main menu and play again buttons"

# match the formatting of the rest of the scenario table (Calibri 11, same as B9:C11)
$newCells = $wsPlay.Range("C12,B13,C13")
$newCells.Font.Name = "Calibri"
$newCells.Font.Size = 11

# --- "SFX Volume" sheet: fix typo "let" -> "left" ---
$wsSfx = $wb.Worksheets.Item("SFX Volume")
$wsSfx.Range("B11").Value = "drag slider left or right"
